$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 15:07:51"
$wsZhCn.Range("H2").Value = "2016-03-21 15:08:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 15:07:55"
$wsDeDe.Range("H2").Value = "2016-03-21 15:08:19"
